# Auto-generated edit script: updates market-price derived columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets to match refreshed
# market data from the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 536.4286
$ws.Range("I2").Value = 713.875
$ws.Range("J2").Value = 299.83334
$ws.Range("K2").Value = 713.875
$ws.Range("L2").Value = 299.83334
$ws.Range("M2").Value = -600.875
$ws.Range("N2").Value = -525.83334
$ws.Range("H41").Value = 731.25
$ws.Range("I41").Value = 953
$ws.Range("J41").Value = 361.66666
$ws.Range("K41").Value = 953
$ws.Range("L41").Value = 361.66666
$ws.Range("M41").Value = -513
$ws.Range("N41").Value = -1241.66666
$ws.Range("H74").Value = 4760.875
$ws.Range("I74").Value = 4370.5713
$ws.Range("J74").Value = 7493
$ws.Range("K74").Value = 4370.5713
$ws.Range("L74").Value = 7493
$ws.Range("M74").Value = -3434.5713
$ws.Range("N74").Value = -9365
$ws.Range("H77").Value = 4760.875
$ws.Range("I77").Value = 4370.5713
$ws.Range("J77").Value = 7493
$ws.Range("K77").Value = 21852.8565
$ws.Range("L77").Value = 37465
$ws.Range("M77").Value = -17172.8565
$ws.Range("N77").Value = -46825
$ws.Range("H116").Value = 473941.12
$ws.Range("I116").Value = 786736.2
$ws.Range("K116").Value = 786736.2
$ws.Range("M116").Value = -783294.2
$ws.Range("H133").Value = 96850.5
$ws.Range("J133").Value = 96850.5
$ws.Range("L133").Value = 96850.5
$ws.Range("N133").Value = -106970.5
$ws.Range("H137").Value = 260463.58
$ws.Range("I137").Value = 458516.56
$ws.Range("J137").Value = 7395.8887
$ws.Range("K137").Value = 1375549.68
$ws.Range("L137").Value = 22187.6661
$ws.Range("M137").Value = -1372999.68
$ws.Range("N137").Value = -27287.6661

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2448.7173
$ws.Range("I32").Value = 2584.3096
$ws.Range("K32").Value = 2584.3096
$ws.Range("M32").Value = -2297.3096
$ws.Range("H45").Value = 10925.947
$ws.Range("I45").Value = 11984.286
$ws.Range("K45").Value = 11984.286
$ws.Range("M45").Value = -11607.286
$ws.Range("H61").Value = 7240.1816
$ws.Range("J61").Value = 5670.6665
$ws.Range("L61").Value = 5670.6665
$ws.Range("N61").Value = -6094.6665
$ws.Range("H63").Value = 1425.625
$ws.Range("I63").Value = 1425.625
$ws.Range("K63").Value = 1425.625
$ws.Range("M63").Value = -739.625
$ws.Range("H66").Value = 1425.625
$ws.Range("I66").Value = 1425.625
$ws.Range("K66").Value = 7128.125
$ws.Range("M66").Value = -3696.125
$ws.Range("H132").Value = 2897.5386
$ws.Range("I132").Value = 2335.9678
$ws.Range("J132").Value = 5073.625
$ws.Range("K132").Value = 7007.903399999999
$ws.Range("L132").Value = 15220.875
$ws.Range("M132").Value = -4477.903399999999
$ws.Range("N132").Value = -20280.875
$ws.Range("H136").Value = 7240.1816
$ws.Range("J136").Value = 5670.6665
$ws.Range("L136").Value = 17011.9995
$ws.Range("N136").Value = -22111.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2739.3845
$ws.Range("I20").Value = 1735.3334
$ws.Range("K20").Value = 1735.3334
$ws.Range("M20").Value = -1488.3334
$ws.Range("H86").Value = 3983.0605
$ws.Range("I86").Value = 4602.6523
$ws.Range("K86").Value = 4602.6523
$ws.Range("M86").Value = -3479.6523
$ws.Range("H89").Value = 3983.0605
$ws.Range("I89").Value = 4602.6523
$ws.Range("K89").Value = 23013.2615
$ws.Range("M89").Value = -17397.2615
$ws.Range("H138").Value = 114598
$ws.Range("J138").Value = 114598
$ws.Range("L138").Value = 114598
$ws.Range("N138").Value = -124878

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 199.5
$ws.Range("I7").Value = 180.33333
$ws.Range("J7").Value = 234
$ws.Range("K7").Value = 180.33333
$ws.Range("L7").Value = 234
$ws.Range("M7").Value = -67.33332999999999
$ws.Range("N7").Value = -460
$ws.Range("H58").Value = 2971.1936
$ws.Range("I58").Value = 1775.0667
$ws.Range("J58").Value = 4092.5625
$ws.Range("K58").Value = 1775.0667
$ws.Range("L58").Value = 4092.5625
$ws.Range("M58").Value = -1572.0667
$ws.Range("N58").Value = -4498.5625
$ws.Range("H62").Value = 209333.33
$ws.Range("J62").Value = 209333.33
$ws.Range("L62").Value = 209333.33
$ws.Range("N62").Value = -210581.33
$ws.Range("H65").Value = 209333.33
$ws.Range("J65").Value = 209333.33
$ws.Range("L65").Value = 1046666.65
$ws.Range("N65").Value = -1052906.65
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("M82:N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("M85:N85").ClearContents()
$ws.Range("H94").Value = 3049.5833
$ws.Range("J94").Value = 1613.1111
$ws.Range("L94").Value = 1613.1111
$ws.Range("N94").Value = -2515.1111
$ws.Range("H132").Value = 24778.041
$ws.Range("I132").Value = 8384.429
$ws.Range("J132").Value = 139533.33
$ws.Range("K132").Value = 25153.287
$ws.Range("L132").Value = 418599.99
$ws.Range("M132").Value = -22623.287
$ws.Range("N132").Value = -423659.99
$ws.Range("H134").Value = 4176733
$ws.Range("I134").Value = 4818230.5
$ws.Range("J134").Value = 6999.5
$ws.Range("K134").Value = 14454691.5
$ws.Range("L134").Value = 20998.5
$ws.Range("M134").Value = -14452156.5
$ws.Range("N134").Value = -26068.5
$ws.Range("H136").Value = 2971.1936
$ws.Range("I136").Value = 1775.0667
$ws.Range("J136").Value = 4092.5625
$ws.Range("K136").Value = 5325.2001
$ws.Range("L136").Value = 12277.6875
$ws.Range("M136").Value = -2775.2001
$ws.Range("N136").Value = -17377.6875
$ws.Range("H137").Value = 66399.39999999999
$ws.Range("J137").Value = 66399.39999999999
$ws.Range("L137").Value = 66399.39999999999
$ws.Range("N137").Value = -76599.39999999999
$ws.Range("H138").Value = 145000
$ws.Range("J138").Value = 145000
$ws.Range("L138").Value = 145000
$ws.Range("N138").Value = -155280
$ws.Range("H140").Value = 67326.664
$ws.Range("J140").Value = 86980
$ws.Range("L140").Value = 86980
$ws.Range("N140").Value = -97340

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 811.7778
$ws.Range("I50").Value = 363.6
$ws.Range("J50").Value = 1372
$ws.Range("K50").Value = 1090.8
$ws.Range("L50").Value = 4116
$ws.Range("M50").Value = -609.8000000000002
$ws.Range("N50").Value = -5078
$ws.Range("H53").Value = 811.7778
$ws.Range("I53").Value = 363.6
$ws.Range("J53").Value = 1372
$ws.Range("K53").Value = 1090.8
$ws.Range("L53").Value = 4116
$ws.Range("M53").Value = -609.8000000000002
$ws.Range("N53").Value = -5078
$ws.Range("H56").Value = 5970.1816
$ws.Range("I56").Value = 5970.1816
$ws.Range("K56").Value = 5970.1816
$ws.Range("M56").Value = -5440.1816

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H132").Value = 3478.182
$ws.Range("I132").Value = 2415.5334
$ws.Range("K132").Value = 7246.600199999999
$ws.Range("M132").Value = -4716.600199999999
$ws.Range("H135").Value = 93554
$ws.Range("J135").Value = 93554
$ws.Range("L135").Value = 93554
$ws.Range("N135").Value = -103694

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3932.7144
$ws.Range("I122").Value = 3116.1667
$ws.Range("K122").Value = 9348.500100000001
$ws.Range("M122").Value = -6898.500100000001
$ws.Range("H132").Value = 1377448.5
$ws.Range("I132").Value = 1898422.4
$ws.Range("K132").Value = 5695267.199999999
$ws.Range("M132").Value = -5692737.199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H128").Value = 37715
$ws.Range("J128").Value = 37715
$ws.Range("L128").Value = 37715
$ws.Range("N128").Value = -47675
$ws.Range("H132").Value = 14595.781
$ws.Range("I132").Value = 22587.611
$ws.Range("K132").Value = 67762.833
$ws.Range("M132").Value = -65232.833
